$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '42.804.29'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = "'" + '2.547.58'
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = "'" + '310.22'
$ws.Range("E5").Value = '  -2.12%  '

$ws.Range("D6").Value = "'" + '98.57'
$ws.Range("E6").Value = '  +0.96%  '

$ws.Range("D7").Value = "'" + '0.571'
$ws.Range("E7").Value = '  -0.61%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  -0.95%  '

$ws.Range("D10").Value = "'" + '35.75'
$ws.Range("E10").Value = '  -1.60%  '

$ws.Range("E11").Value = '  -1.61%  '

$ws.Range("D12").Value = "'" + '7.41'
$ws.Range("E12").Value = '  -3.06%  '

$ws.Range("E13").Value = '  -1.56%  '

$ws.Range("D14").Value = "'" + '2.940.15'
$ws.Range("E14").Value = '  +0.35%  '

$ws.Range("D15").Value = "'" + '15.95'
$ws.Range("E15").Value = '  +4.89%  '

$ws.Range("D16").Value = "'" + '2.577.03'
$ws.Range("E16").Value = '  +1.21%  '

$ws.Range("E17").Value = '  -1.48%  '

$ws.Range("D18").Value = "'" + '42.831.64'
$ws.Range("E18").Value = '  -0.38%  '

$ws.Range("D19").Value = "'" + '6.75'

$ws.Range("D20").Value = "'" + '12.42'
$ws.Range("E20").Value = '  -3.21%  '

$ws.Range("E21").Value = '  -1.27%  '

$ws.Range("D22").Value = "'" + '69.35'
$ws.Range("E22").Value = '  -0.89%  '

$ws.Range("D23").Value = "'" + '248.15'
$ws.Range("E23").Value = '  -2.58%  '

$ws.Range("D24").Value = "'" + '2.93'
$ws.Range("E24").Value = '  -0.74%  '

$ws.Range("E25").Value = '  +0.27%  '

$ws.Range("D26").Value = "'" + '26.67'
$ws.Range("E26").Value = '  +0.18%  '

$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("D28").Value = "'" + '2.36'
$ws.Range("E28").Value = '  -2.14%  '

$ws.Range("D29").Value = "'" + '40.25'
$ws.Range("E29").Value = '  -2.00%  '

$ws.Range("D30").Value = "'" + '10.17'
$ws.Range("E30").Value = '  -2.91%  '

$ws.Range("D31").Value = "'" + '158.32'
$ws.Range("E31").Value = '  -0.04%  '

$ws.Range("E32").Value = '  -3.00%  '

$ws.Range("E33").Value = '  +0.83%  '

$ws.Range("D34").Value = "'" + '3.30'
$ws.Range("E34").Value = '  -2.07%  '

$ws.Range("E35").Value = '  -3.92%  '

$ws.Range("E36").Value = '  -3.55%  '

$ws.Range("D37").Value = "'" + '2.61'
$ws.Range("E37").Value = '  +4.63%  '

$ws.Range("D38").Value = "'" + '18.31'
$ws.Range("E38").Value = '  -3.82%  '

$ws.Range("E39").Value = '  -1.52%  '

$ws.Range("E40").Value = '  -0.72%  '

$ws.Range("D41").Value = "'" + '22.46'
$ws.Range("E41").Value = '  +2.44%  '

$ws.Range("E42").Value = '  +6.48%  '

$ws.Range("E43").Value = '  -0.10%  '

$ws.Range("E44").Value = '  -1.43%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = "'" + '3.21'
$ws.Range("E45").Value = '  -2.70%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = "'" + '1.996.78'
$ws.Range("E46").Value = '  -1.45%  '

$ws.Range("D47").Value = "'" + '9.09'
$ws.Range("E47").Value = '  -0.24%  '

$ws.Range("D48").Value = "'" + '2.787.27'
$ws.Range("E48").Value = '  +0.09%  '

$ws.Range("D49").Value = "'" + '81.32'
$ws.Range("E49").Value = '  -3.93%  '

$ws.Range("D50").Value = "'" + '0.192'
$ws.Range("E50").Value = '  -0.07%  '

$ws.Range("D51").Value = "'" + '73.61'
$ws.Range("E51").Value = '  -4.55%  '
